$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Users sheet: swap out a deal-team member name
#    (Liz Hedgcock -> Blaise Brunda)
# ------------------------------------------------------------------
$usersSheet = $wb.Worksheets.Item("Users")
$usersSheet.Range("B2").Value = "Blaise Brunda"
$usersSheet.Range("B23").Select()

# ------------------------------------------------------------------
# 2. OppDealTeamMembers: add two new team members at the bottom
# ------------------------------------------------------------------
$oppSheet = $wb.Worksheets.Item("OppDealTeamMembers")
$oppSheet.Range("A28").Select()
$oppSheet.Range("A29").Value = "Tom Seward"
$oppSheet.Range("A30").Value = "Lucy Gao"
$oppSheet.Application.ActiveWindow.ScrollRow = 5
$oppSheet.Range("A28").Select()

# ------------------------------------------------------------------
# 3. EngDealTeamMembers: swap out a deal-team member name
#    (Timothy Kang -> Zev Litwin)
# ------------------------------------------------------------------
$engSheet = $wb.Worksheets.Item("EngDealTeamMembers")
$engSheet.Range("A2").Value = "Zev Litwin"
$engSheet.Range("A4").Select()

# ------------------------------------------------------------------
# 4. Reorder tabs: move "Users" to the front of the workbook
#    (before "AddOpportunity")
# ------------------------------------------------------------------
$usersSheet.Move($wb.Worksheets.Item(1))

# ------------------------------------------------------------------
# 5. Restore the previously-active sheet/selection
#    (OppDealTeamMembers stays the selected/active tab)
# ------------------------------------------------------------------
$oppSheet.Select()
$oppSheet.Range("A28").Select()
